$d = $word.ActiveDocument

$replacements = @(
    @("2024-12-22 Sunday", "2024-12-23 Monday"),
    @("97×86=", "74×58="),
    @("29×94=", "18×47="),
    @("31×39=", "78×64="),
    @("12×63=", "68×30="),
    @("19×17=", "94×29="),
    @("73×48=", "21×50="),
    @("16×87=", "73×87="),
    @("91×38=", "33×49="),
    @("16×88=", "66×91="),
    @("32×43=", "62×15="),
    @("47×45=", "15×56="),
    @("14×75=", "15×69="),
    @("56×97=", "90×99="),
    @("38×24=", "56×65="),
    @("77×84=", "81×11="),
    @("20×35=", "58×63="),
    @("84×69=", "52×41="),
    @("54×95=", "55×38="),
    @("23×47=", "13×14="),
    @("57×85=", "40×27="),
    @("62×21=", "83×62="),
    @("34×52=", "78×83="),
    @("63×13=", "71×46="),
    @("73×31=", "55×15="),
    @("62×34=", "51×69=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
